$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("D2").Value = 2
$ws.Range("C3").Value = 1

$ws.Range("D2").Select()
